$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'283.21"
$ws.Range("E2").Value = "'1.79%"
$ws.Range("D3").Value = "'28.34"
$ws.Range("E3").Value = "'4.12%"
$ws.Range("D4").Value = "'5.030"
$ws.Range("E4").Value = "'3.26%"
$ws.Range("D5").Value = "'0.06502"
$ws.Range("E5").Value = "'1.18%"
$ws.Range("D6").Value = "'7.220"
$ws.Range("E6").Value = "'3.07%"
$ws.Range("D7").Value = "'1.385"
$ws.Range("E7").Value = "'14.71%"
$ws.Range("D8").Value = "'0.9181"
$ws.Range("E8").Value = "'3.49%"
$ws.Range("D9").Value = "'0.1547"
$ws.Range("E9").Value = "'-0.01%"
$ws.Range("D10").Value = "'0.06369"
$ws.Range("E10").Value = "'24.71%"
$ws.Range("D11").Value = "'0.07585"
$ws.Range("E11").Value = "'1.05%"
$ws.Range("D12").Value = "'0.02850"
$ws.Range("E12").Value = "'-1.18%"
$ws.Range("D13").Value = "'0.08989"
$ws.Range("E13").Value = "'0.21%"
$ws.Range("D14").Value = "'0.001585"
$ws.Range("E14").Value = "'0.74%"
$ws.Range("D15").Value = "'0.0006363"
$ws.Range("E15").Value = "'-0.50%"
$ws.Range("D16").Value = "'0.006119"
$ws.Range("E16").Value = "'0.76%"
$ws.Range("D17").Value = "'3.443"
$ws.Range("E17").Value = "'-0.96%"
$ws.Range("E18").Value = "'1.66%"
$ws.Range("E19").Value = "'1.08%"
$ws.Range("E20").Value = "'-0.06%"
$ws.Range("E21").Value = "'-0.92%"
$ws.Range("D22").Value = "'3.981"
$ws.Range("E22").Value = "'1.89%"
$ws.Range("E24").Value = "'0.44%"
$ws.Range("D25").Value = "'0.001181"
$ws.Range("E25").Value = "'0.36%"
$ws.Range("D26").Value = "'0.004455"
$ws.Range("E26").Value = "'14.85%"
$ws.Range("E28").Value = "'1.64%"
$ws.Range("D29").Value = "'0.0001617"
$ws.Range("E29").Value = "'-1.61%"
$ws.Range("D40").Value = "'0.04115"
$ws.Range("E40").Value = "'-0.20%"
$ws.Range("D41").Value = "'0.006651"
$ws.Range("E41").Value = "'-1.73%"
$ws.Range("E42").Value = "'4.91%"
$ws.Range("D43").Value = "'0.002139"
$ws.Range("E43").Value = "'11.39%"
$ws.Range("D44").Value = "'0.01155"
$ws.Range("E44").Value = "'-2.51%"
$ws.Range("D45").Value = "'0.00005613"
$ws.Range("E45").Value = "'5.18%"
$ws.Range("D46").Value = "'1.978"
$ws.Range("E46").Value = "'17.33%"
$ws.Range("E47").Value = "'-0.06%"
